$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = "control system"
$ws.Range("B4").Value = "27/06/2025, 13:31"
$ws.Range("C4").Value = "WIVGH"
$ws.Range("D4").Value = "'2025-06-27"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'25"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "Mecânico"
$ws.Range("G4").Value = "Castanha/Cardã"
$ws.Range("H4").Value = "'"
$ws.Range("H4").Style = "Normal"

# Row 5
$ws.Range("A5").Value = "control system"
$ws.Range("B5").Value = "27/06/2025, 13:34"
$ws.Range("C5").Value = "FKB58"
$ws.Range("D5").Value = "'2025-06-27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'25"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "Elétrico"
$ws.Range("G5").Value = "Anel Coletor"
$ws.Range("H5").Value = "'"
$ws.Range("H5").Style = "Normal"
